$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.830.05"
$ws.Range("D3").Value = "2.275.90"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.44%  "
$ws.Range("E6").Value = "  +5.77%  "
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0801"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "2.628.26"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "2.278.32"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("E18").Value = "  +3.62%  "
$ws.Range("D19").Value = "41.778.14"
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.58%  "
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.02%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0750"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("E37").Value = "  +2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.51%  "
$ws.Range("E39").Value = "  +5.22%  "
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("E42").Value = "  +5.44%  "
$ws.Range("D43").Value = "2.077.39"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0280"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.35%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.71%  "
